$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "23.363.58"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  -1.53%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.629.28"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  -1.78%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.001"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "  +0.05%  "

$ws.Range("E5").Value = "  +0.09%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "298.25"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -1.80%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3755"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  -1.54%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "50.44"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  -1.57%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.3487"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  -3.66%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.08027"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -2.28%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.201"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -3.31%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.001"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  +0.09%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "21.80"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  -4.01%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.307"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -3.41%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.218"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  -2.87%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.00001194"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  -3.26%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.630.86"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  -1.43%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "94.69"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  -2.99%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06935"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -0.90%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.618"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -3.29%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "17.31"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -2.25%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "12.36"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -3.79%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "23.356.47"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -1.54%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.421"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -3.81%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.960"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -3.12%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "20.71"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -2.77%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "152.09"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +0.10%  "

$ws.Range("E29").Value = "  -0.48%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "131.55"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -2.31%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.811.16"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -1.25%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.761"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -3.13%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.117"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -3.38%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "11.10"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -7.29%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.9698"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -8.96%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.02674"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -5.23%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.08711"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -0.83%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.2421"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -4.30%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.843"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -4.78%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.06690"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -5.47%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "12.69"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -3.25%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.6793"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -3.55%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.283"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -4.18%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "15.42"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -3.74%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.0000"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +0.09%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.6305"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -3.79%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.228"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -4.31%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.886"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -2.33%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "126.53"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -1.28%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.07638"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -3.96%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.215"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +1.57%  "
